$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Qminus1)
$ws.Range("B2").Value = -0.076860368503937
$ws.Range("C2").Value = 0.803884330424265
$ws.Range("D2").Value = 1.335363848599442
$ws.Range("E2").Value = 1.155579442790258
$ws.Range("F2").Value = 1.18015413262448
$ws.Range("G2").Value = 22

# Row 3 (Q0)
$ws.Range("B3").Value = 0.08565359546942733
$ws.Range("C3").Value = 0.9749789586808505
$ws.Range("D3").Value = 1.486598319041417
$ws.Range("E3").Value = 1.219261382576114
$ws.Range("F3").Value = 1.243583914527342
$ws.Range("G3").Value = 23

# Row 4 (Q1)
$ws.Range("B4").Value = 0.1093541282614182
$ws.Range("C4").Value = 1.574993429461694
$ws.Range("D4").Value = 10.32836645436715
$ws.Range("E4").Value = 3.213777598771756
$ws.Range("F4").Value = 3.287501431018957
$ws.Range("G4").Value = 22
